$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new quotation row (2025-09-28) as row 24, matching the
# formatting of the previous data row (row 23).
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(24, 1).Value = 45928
$ws.Cells.Item(24, 2).Value = "21,0192"
$ws.Cells.Item(24, 3).Value = "15,0785"
$ws.Cells.Item(24, 4).Value = "14,8412"
$ws.Cells.Item(24, 5).Value = "14,8412"
